$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AG2").Value = "mineral"
$ws.Range("AG3").Value = "no_zero"
$ws.Range("AG4:AG9").Value = "oth"
$ws.Range("AG10").Value = "no_zero"
$ws.Range("AG11:AG13").Value = "mineraloth"
$ws.Range("AG14").Value = "mineral"
$ws.Range("AG15").Value = "mineraloth"
$ws.Range("AG16").Value = "mineral"
$ws.Range("AG17").Value = "mineraloth"
$ws.Range("AG18").Value = "no_zero"
$ws.Range("AG19").Value = "oth"
$ws.Range("AG20").Value = "no_zero"
$ws.Range("AG21").Value = "mineraloth"
$ws.Range("AG22:AG24").Value = "oth"
$ws.Range("AG25").Value = "no_zero"
$ws.Range("AG26").Value = "mineral"
$ws.Range("AG27").Value = "oth"
$ws.Range("AG28").Value = "mineraloth"
$ws.Range("AG29:AG30").Value = "oth"
$ws.Range("AG31:AG32").Value = "no_zero"
$ws.Range("AG33:AG34").Value = "oth"
$ws.Range("AG35").Value = "mineraloth"
$ws.Range("AG36").Value = "mineral"
$ws.Range("AG37").Value = "oth"
$ws.Range("AG38").Value = "mineraloth"
$ws.Range("AG39:AG40").Value = "no_zero"
$ws.Range("AG41").Value = "mineral"
$ws.Range("AG42").Value = "no_zero"
$ws.Range("AG43").Value = "mineraloth"
$ws.Range("AG44:AG45").Value = "oth"
$ws.Range("AG46").Value = "mineraloth"
$ws.Range("AG47").Value = "no_zero"
$ws.Range("AG48").Value = "mineraloth"
$ws.Range("AG49").Value = "mineral"
$ws.Range("AG50").Value = "no_zero"
$ws.Range("AG51").Value = "mineral"
$ws.Range("AG52").Value = "oth"
$ws.Range("AG53").Value = "no_zero"
$ws.Range("AG54:AG55").Value = "oth"
$ws.Range("AG56").Value = "mineraloth"
$ws.Range("AG57:AG58").Value = "oth"
$ws.Range("AG59").Value = "mineraloth"
$ws.Range("AG60").Value = "no_zero"
$ws.Range("AG61").Value = "oth"
$ws.Range("AG62").Value = "mineral"
$ws.Range("AG63").Value = "mineraloth"
$ws.Range("AG64").Value = "no_zero"
$ws.Range("AG65").Value = "oth"
$ws.Range("AG66").Value = "mineraloth"
$ws.Range("AG67").Value = "no_zero"
$ws.Range("AG68:AG70").Value = "mineraloth"
$ws.Range("AG71").Value = "mineral"
$ws.Range("AG72:AG75").Value = "no_zero"
$ws.Range("AG76").Value = "mineral"
$ws.Range("AG77").Value = "no_zero"
$ws.Range("AG78").Value = "mineraloth"
$ws.Range("AG79").Value = "no_zero"
$ws.Range("AG80").Value = "oth"
$ws.Range("AG81").Value = "no_zero"
$ws.Range("AG82").Value = "oth"
$ws.Range("AG83").Value = "no_zero"
$ws.Range("AG84").Value = "mineral"
$ws.Range("AG85:AG86").Value = "oth"
$ws.Range("AG87").Value = "mineraloth"
$ws.Range("AG88:AG90").Value = "oth"
$ws.Range("AG91").Value = "no_zero"
$ws.Range("AG92").Value = "mineraloth"
$ws.Range("AG93:AG95").Value = "oth"
$ws.Range("AG96").Value = "mineraloth"
$ws.Range("AG97:AG99").Value = "no_zero"
$ws.Range("AG100").Value = "oth"
$ws.Range("AG101:AG102").Value = "mineraloth"
$ws.Range("AG103").Value = "oth"
$ws.Range("AG104").Value = "no_zero"
$ws.Range("AG105").Value = "oth"
$ws.Range("AG106:AG107").Value = "no_zero"
$ws.Range("AG108").Value = "oth"
$ws.Range("AG109").Value = "mineraloth"
$ws.Range("AG110").Value = "oth"
$ws.Range("AG111").Value = "mineraloth"
$ws.Range("AG112:AG115").Value = "oth"
$ws.Range("AG116").Value = "mineraloth"
$ws.Range("AG117").Value = "oth"
$ws.Range("AG118").Value = "mineraloth"
$ws.Range("AG119:AG120").Value = "oth"
$ws.Range("AG121").Value = "no_zero"
$ws.Range("AG122:AG123").Value = "oth"
$ws.Range("AG124").Value = "mineraloth"
$ws.Range("AG125:AG126").Value = "oth"
$ws.Range("AG127").Value = "mineral"
$ws.Range("AG128").Value = "oth"
$ws.Range("AG129:AG131").Value = "mineraloth"
$ws.Range("AG132:AG134").Value = "oth"
$ws.Range("AG135:AG136").Value = "mineraloth"
$ws.Range("AG137:AG140").Value = "no_zero"
$ws.Range("AG141").Value = "oth"
$ws.Range("AG142").Value = "no_zero"
$ws.Range("AG143").Value = "oth"
$ws.Range("AG144").Value = "mineraloth"
$ws.Range("AG145:AG147").Value = "oth"
$ws.Range("AG148").Value = "mineraloth"
$ws.Range("AG149").Value = "oth"
$ws.Range("AG150").Value = "mineraloth"
$ws.Range("AG151:AG155").Value = "no_zero"
$ws.Range("AG156").Value = "oth"
$ws.Range("AG157").Value = "no_zero"
$ws.Range("AG158").Value = "mineral"
$ws.Range("AG159").Value = "no_zero"
$ws.Range("AG160").Value = "mineraloth"
$ws.Range("AG161").Value = "no_zero"
$ws.Range("AG162").Value = "oth"
$ws.Range("AG163").Value = "no_zero"
$ws.Range("AG164").Value = "mineraloth"
$ws.Range("AG165").Value = "no_zero"
$ws.Range("AG166").Value = "oth"
$ws.Range("AG167").Value = "no_zero"
$ws.Range("AG168").Value = "mineral"
$ws.Range("AG169:AG170").Value = "no_zero"
$ws.Range("AG171").Value = "oth"
$ws.Range("AG172").Value = "mineral"
$ws.Range("AG173").Value = "mineraloth"
$ws.Range("AG174").Value = "no_zero"
$ws.Range("AG175").Value = "oth"
$ws.Range("AG176").Value = "mineraloth"
$ws.Range("AG177:AG179").Value = "oth"
$ws.Range("AG180").Value = "mineraloth"
